# "Generalized SteppedController to diff types of selection"
#
# selection-control sheet:
#   - C2: "non-chemical" -> "chemical"
#   - G2:G17: 7 -> 3
#   - H2:H17: 6 -> 2
#   - I2:I17: 0.1 -> 0.3
#   - J2:J17: 0.12 -> 0.4
#   - New header comment on E1
#
# selection-step_generation sheet:
#   - D3: 50 -> 20

$wb = $excel.ActiveWorkbook

$wsControl = $wb.Worksheets.Item("selection-control")
$wsStepGen = $wb.Worksheets.Item("selection-step_generation")

# --- selection-control sheet -------------------------------------------------

$wsControl.Range("C2").Value = "chemical"

$wsControl.Range("G2:G17").Value = 3
$wsControl.Range("H2:H17").Value = 2
$wsControl.Range("I2:I17").Value = 0.3
$wsControl.Range("J2:J17").Value = 0.4

$wsControl.Range("E1").AddComment("number of growth curves required at start of experiment before beginning selection") | Out-Null

# --- selection-step_generation sheet -----------------------------------------

$wsStepGen.Range("D3").Value = 20

# --- selections / active cells -----------------------------------------------

$wsStepGen.Activate()
$wsStepGen.Range("D3").Select()

$wsControl.Activate()
$wsControl.Range("E8").Select()
